# Generate Report for Handback
# Adds a new tracked file (98b07430-8a3e-41f3-b911-bb49260c42c3.md) as row 4
# to the "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$fileId   = "98b07430-8a3e-41f3-b911-bb49260c42c3"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$genericCommit = "8fe89bfe0424f83998518d355cdbf4f946c64953"
$xlfCommit     = "0f8894d879860b5e898e7ea0202adf798ca825da"

function Set-TextValue($range, $text) {
    # Forces a literal text value (prevents Excel from coercing look-alike
    # booleans/dates/numbers), then strips the resulting quote-prefix flag
    # so the cell style falls back to the sheet's normal style.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-TextValue $wsOverview.Range("A4") $mdName
Set-TextValue $wsOverview.Range("B4") $mdPath
Set-TextValue $wsOverview.Range("C4") ".md"
Set-TextValue $wsOverview.Range("E4") "Handed back: in sync with en-US"
Set-TextValue $wsOverview.Range("F4") "Handed back: in sync with en-US"
Set-TextValue $wsOverview.Range("G4") "2016-10-18 05:01:21"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$genericCommit/e2e/$mdName",
    "",
    "",
    $mdPath
) | Out-Null

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-TextValue $wsZhCn.Range("A4") $mdName
Set-TextValue $wsZhCn.Range("B4") ".md"
Set-TextValue $wsZhCn.Range("C4") "Handed back: in sync with en-US"
Set-TextValue $wsZhCn.Range("D4") "e2e"
Set-TextValue $wsZhCn.Range("E4") "ht"
Set-TextValue $wsZhCn.Range("F4") "True"
Set-TextValue $wsZhCn.Range("G4") "$fileId.$xlfCommit.zh-cn.xlf"
Set-TextValue $wsZhCn.Range("H4") "2016-10-18 05:01:06"
Set-TextValue $wsZhCn.Range("I4") $mdName
Set-TextValue $wsZhCn.Range("J4") "$fileId.$xlfCommit.zh-cn.xlf"
Set-TextValue $wsZhCn.Range("K4") "2016-10-18 05:01:59"
Set-TextValue $wsZhCn.Range("L4") ""
Set-TextValue $wsZhCn.Range("M4") "True"
Set-TextValue $wsZhCn.Range("N4") ""
Set-TextValue $wsZhCn.Range("O4") "False"
Set-TextValue $wsZhCn.Range("P4") ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$genericCommit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$xlfCommit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-TextValue $wsDeDe.Range("A4") $mdName
Set-TextValue $wsDeDe.Range("B4") ".md"
Set-TextValue $wsDeDe.Range("C4") "Handed back: in sync with en-US"
Set-TextValue $wsDeDe.Range("D4") "e2e"
Set-TextValue $wsDeDe.Range("E4") "ht"
Set-TextValue $wsDeDe.Range("F4") "True"
Set-TextValue $wsDeDe.Range("G4") "$fileId.$xlfCommit.de-de.xlf"
Set-TextValue $wsDeDe.Range("H4") "2016-10-18 05:01:21"
Set-TextValue $wsDeDe.Range("I4") $mdName
Set-TextValue $wsDeDe.Range("J4") "$fileId.$xlfCommit.de-de.xlf"
Set-TextValue $wsDeDe.Range("K4") "2016-10-18 05:02:23"
Set-TextValue $wsDeDe.Range("L4") ""
Set-TextValue $wsDeDe.Range("M4") "True"
Set-TextValue $wsDeDe.Range("N4") ""
Set-TextValue $wsDeDe.Range("O4") "False"
Set-TextValue $wsDeDe.Range("P4") ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$genericCommit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$xlfCommit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Host "Report row added for $mdName"
